# Auto-generated Excel COM-interop script
# Applies scheduled market-price / profit updates to the Jenova_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 210.2
$ws.Cells.Item(4, 9).Value = 210.2
$ws.Cells.Item(4, 11).Value = 210.2
$ws.Cells.Item(4, 13).Value = -96.19999999999999

$ws.Cells.Item(40, 8).Value = 7384.643
$ws.Cells.Item(40, 9).Value = 3599.8
$ws.Cells.Item(40, 10).Value = 9487.333000000001
$ws.Cells.Item(40, 11).Value = 3599.8
$ws.Cells.Item(40, 12).Value = 9487.333000000001
$ws.Cells.Item(40, 13).Value = -3424.8
$ws.Cells.Item(40, 14).Value = -9837.333000000001

$ws.Cells.Item(47, 8).Value = 11987
$ws.Cells.Item(47, 9).Value = 9000
$ws.Cells.Item(47, 11).Value = 9000
$ws.Cells.Item(47, 13).Value = -8028

$ws.Cells.Item(100, 8).Value = 5576
$ws.Cells.Item(100, 9).Value = 1826.7693
$ws.Cells.Item(100, 11).Value = 1826.7693
$ws.Cells.Item(100, 13).Value = -1285.7693

$ws.Cells.Item(106, 8).Value = 2498.4783
$ws.Cells.Item(106, 9).Value = 2524.8667
$ws.Cells.Item(106, 10).Value = 2449
$ws.Cells.Item(106, 11).Value = 2524.8667
$ws.Cells.Item(106, 12).Value = 2449
$ws.Cells.Item(106, 13).Value = -1893.8667
$ws.Cells.Item(106, 14).Value = -3711

$ws.Cells.Item(132, 8).Value = 2286.1035
$ws.Cells.Item(132, 9).Value = 1091.8182
$ws.Cells.Item(132, 11).Value = 3275.4546
$ws.Cells.Item(132, 13).Value = -745.4546

$ws.Cells.Item(137, 8).Value = 5017.353
$ws.Cells.Item(137, 9).Value = 3022.4285
$ws.Cells.Item(137, 10).Value = 8239.923000000001
$ws.Cells.Item(137, 11).Value = 9067.2855
$ws.Cells.Item(137, 12).Value = 24719.769
$ws.Cells.Item(137, 13).Value = -6517.2855
$ws.Cells.Item(137, 14).Value = -29819.769

$ws.Cells.Item(138, 8).Value = 5941.304
$ws.Cells.Item(138, 9).Value = 5671
$ws.Cells.Item(138, 10).Value = 6026.2573
$ws.Cells.Item(138, 11).Value = 17013
$ws.Cells.Item(138, 12).Value = 18078.7719
$ws.Cells.Item(138, 13).Value = -11873
$ws.Cells.Item(138, 14).Value = -28358.7719

$ws.Cells.Item(141, 8).Value = 3963.4
$ws.Cells.Item(141, 9).Value = 1772.3334
$ws.Cells.Item(141, 10).Value = 7250
$ws.Cells.Item(141, 11).Value = 5317.0002
$ws.Cells.Item(141, 12).Value = 21750
$ws.Cells.Item(141, 13).Value = -137.0002000000004
$ws.Cells.Item(141, 14).Value = -32110

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 854
$ws.Cells.Item(19, 9).Value = 854
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 854
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -625
$ws.Cells.Item(19, 14).ClearContents()

$ws.Cells.Item(57, 8).Value = 10000
$ws.Cells.Item(57, 9).Value = 10000
$ws.Cells.Item(57, 11).Value = 10000
$ws.Cells.Item(57, 13).Value = -9516

$ws.Cells.Item(74, 8).Value = 4657.091
$ws.Cells.Item(74, 9).Value = 3028.625
$ws.Cells.Item(74, 11).Value = 3028.625
$ws.Cells.Item(74, 13).Value = -2154.625

$ws.Cells.Item(77, 8).Value = 4657.091
$ws.Cells.Item(77, 9).Value = 3028.625
$ws.Cells.Item(77, 11).Value = 15143.125
$ws.Cells.Item(77, 13).Value = -10775.125

$ws.Cells.Item(102, 8).Value = 1865.5
$ws.Cells.Item(102, 9).Value = 1865.5
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1865.5
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -243.5
$ws.Cells.Item(102, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 5032.778
$ws.Cells.Item(122, 9).Value = 3304.7778
$ws.Cells.Item(122, 10).Value = 6760.778
$ws.Cells.Item(122, 11).Value = 9914.3334
$ws.Cells.Item(122, 12).Value = 20282.334
$ws.Cells.Item(122, 13).Value = -7464.3334
$ws.Cells.Item(122, 14).Value = -25182.334

$ws.Cells.Item(132, 8).Value = 4044.6182
$ws.Cells.Item(132, 9).Value = 1296.6428
$ws.Cells.Item(132, 10).Value = 12922.692
$ws.Cells.Item(132, 11).Value = 3889.9284
$ws.Cells.Item(132, 12).Value = 38768.076
$ws.Cells.Item(132, 13).Value = -1359.9284
$ws.Cells.Item(132, 14).Value = -43828.076

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(102, 8).Value = 5542
$ws.Cells.Item(102, 9).Value = 5542
$ws.Cells.Item(102, 11).Value = 5542
$ws.Cells.Item(102, 13).Value = -2297

$ws.Cells.Item(134, 8).Value = 2490.718
$ws.Cells.Item(134, 9).Value = 1732.9697
$ws.Cells.Item(134, 11).Value = 5198.909100000001
$ws.Cells.Item(134, 13).Value = -2663.909100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4000
$ws.Cells.Item(31, 9).Value = 1998.25
$ws.Cells.Item(31, 10).Value = 6669
$ws.Cells.Item(31, 11).Value = 1998.25
$ws.Cells.Item(31, 12).Value = 6669
$ws.Cells.Item(31, 13).Value = -1703.25
$ws.Cells.Item(31, 14).Value = -7259

$ws.Cells.Item(34, 8).Value = 4000
$ws.Cells.Item(34, 9).Value = 1998.25
$ws.Cells.Item(34, 10).Value = 6669
$ws.Cells.Item(34, 11).Value = 1998.25
$ws.Cells.Item(34, 12).Value = 6669
$ws.Cells.Item(34, 13).Value = -1796.25
$ws.Cells.Item(34, 14).Value = -7073

$ws.Cells.Item(122, 8).Value = 2801.5
$ws.Cells.Item(122, 9).Value = 1399.625
$ws.Cells.Item(122, 10).Value = 4670.6665
$ws.Cells.Item(122, 11).Value = 4198.875
$ws.Cells.Item(122, 12).Value = 14011.9995
$ws.Cells.Item(122, 13).Value = -1748.875
$ws.Cells.Item(122, 14).Value = -18911.9995

$ws.Cells.Item(132, 8).Value = 6108.1665
$ws.Cells.Item(132, 9).Value = 5277.857
$ws.Cells.Item(132, 11).Value = 15833.571
$ws.Cells.Item(132, 13).Value = -13303.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 51378.688
$ws.Cells.Item(5, 9).Value = 81040.10000000001
$ws.Cells.Item(5, 10).Value = 1943
$ws.Cells.Item(5, 11).Value = 243120.3
$ws.Cells.Item(5, 12).Value = 5829
$ws.Cells.Item(5, 13).Value = -243008.3
$ws.Cells.Item(5, 14).Value = -6053

$ws.Cells.Item(33, 8).Value = 2849022.5
$ws.Cells.Item(33, 9).Value = 3086440
$ws.Cells.Item(33, 10).Value = 14
$ws.Cells.Item(33, 11).Value = 18518640
$ws.Cells.Item(33, 12).Value = 84
$ws.Cells.Item(33, 13).Value = -18518357
$ws.Cells.Item(33, 14).Value = -650

$ws.Cells.Item(70, 8).Value = 5650
$ws.Cells.Item(70, 10).Value = 14000
$ws.Cells.Item(70, 12).Value = 42000
$ws.Cells.Item(70, 14).Value = -42630

$ws.Cells.Item(73, 8).Value = 5650
$ws.Cells.Item(73, 10).Value = 14000
$ws.Cells.Item(73, 12).Value = 42000
$ws.Cells.Item(73, 14).Value = -44184

$ws.Cells.Item(114, 8).Value = 1470.125
$ws.Cells.Item(114, 9).Value = 702.25
$ws.Cells.Item(114, 10).Value = 2238
$ws.Cells.Item(114, 11).Value = 2106.75
$ws.Cells.Item(114, 12).Value = 6714
$ws.Cells.Item(114, 13).Value = 1147.25
$ws.Cells.Item(114, 14).Value = -13222

$ws.Cells.Item(119, 8).Value = 1574
$ws.Cells.Item(119, 9).Value = 943.8
$ws.Cells.Item(119, 11).Value = 2831.4
$ws.Cells.Item(119, 13).Value = 2006.6

$ws.Cells.Item(129, 8).Value = 2588.7273
$ws.Cells.Item(129, 9).Value = 652.75
$ws.Cells.Item(129, 10).Value = 3695
$ws.Cells.Item(129, 11).Value = 1958.25
$ws.Cells.Item(129, 12).Value = 11085
$ws.Cells.Item(129, 13).Value = 3041.75
$ws.Cells.Item(129, 14).Value = -21085

$ws.Cells.Item(131, 8).Value = 4689.28
$ws.Cells.Item(131, 10).Value = 8814.75
$ws.Cells.Item(131, 12).Value = 26444.25
$ws.Cells.Item(131, 14).Value = -36524.25

$ws.Cells.Item(135, 8).Value = 51378.688
$ws.Cells.Item(135, 9).Value = 81040.10000000001
$ws.Cells.Item(135, 10).Value = 1943
$ws.Cells.Item(135, 11).Value = 729360.9
$ws.Cells.Item(135, 12).Value = 17487
$ws.Cells.Item(135, 13).Value = -726825.9
$ws.Cells.Item(135, 14).Value = -22557

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(42, 8).Value = 68000
$ws.Cells.Item(42, 10).Value = 68000
$ws.Cells.Item(42, 12).Value = 68000
$ws.Cells.Item(42, 14).Value = -68970

$ws.Cells.Item(102, 8).Value = 2234
$ws.Cells.Item(102, 9).Value = 2408.3333
$ws.Cells.Item(102, 10).Value = 1449.5
$ws.Cells.Item(102, 11).Value = 2408.3333
$ws.Cells.Item(102, 12).Value = 1449.5
$ws.Cells.Item(102, 13).Value = -786.3332999999998
$ws.Cells.Item(102, 14).Value = -4693.5

$ws.Cells.Item(113, 8).Value = 473839.2
$ws.Cells.Item(113, 10).Value = 34927
$ws.Cells.Item(113, 12).Value = 34927
$ws.Cells.Item(113, 14).Value = -39267

$ws.Cells.Item(115, 8).Value = 68000
$ws.Cells.Item(115, 10).Value = 68000
$ws.Cells.Item(115, 12).Value = 68000
$ws.Cells.Item(115, 14).Value = -70350

$ws.Cells.Item(126, 8).Value = 111114340
$ws.Cells.Item(126, 9).Value = 200002940
$ws.Cells.Item(126, 10).Value = 3577.75
$ws.Cells.Item(126, 11).Value = 600008820
$ws.Cells.Item(126, 12).Value = 10733.25
$ws.Cells.Item(126, 13).Value = -600006350
$ws.Cells.Item(126, 14).Value = -15673.25

$ws.Cells.Item(132, 8).Value = 259756.72
$ws.Cells.Item(132, 9).Value = 288760.4
$ws.Cells.Item(132, 11).Value = 866281.2000000001
$ws.Cells.Item(132, 13).Value = -863751.2000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 521.25
$ws.Cells.Item(16, 9).Value = 526.5454999999999
$ws.Cells.Item(16, 10).Value = 463
$ws.Cells.Item(16, 11).Value = 526.5454999999999
$ws.Cells.Item(16, 12).Value = 463
$ws.Cells.Item(16, 13).Value = -356.5454999999999
$ws.Cells.Item(16, 14).Value = -803

$ws.Cells.Item(38, 8).Value = 13033
$ws.Cells.Item(38, 10).Value = 13033
$ws.Cells.Item(38, 12).Value = 13033
$ws.Cells.Item(38, 14).Value = -13853

$ws.Cells.Item(40, 8).Value = 438728.53
$ws.Cells.Item(40, 9).Value = 627877.3
$ws.Cells.Item(40, 10).Value = 6388.4287
$ws.Cells.Item(40, 11).Value = 627877.3
$ws.Cells.Item(40, 12).Value = 6388.4287
$ws.Cells.Item(40, 13).Value = -627741.3
$ws.Cells.Item(40, 14).Value = -6660.4287

$ws.Cells.Item(61, 8).Value = 4286.448
$ws.Cells.Item(61, 9).Value = 2832.9375
$ws.Cells.Item(61, 11).Value = 2832.9375
$ws.Cells.Item(61, 13).Value = -2630.9375

$ws.Cells.Item(100, 8).Value = 13428.571
$ws.Cells.Item(100, 10).Value = 70000
$ws.Cells.Item(100, 12).Value = 70000
$ws.Cells.Item(100, 14).Value = -71082

$ws.Cells.Item(113, 8).Value = 4286.448
$ws.Cells.Item(113, 9).Value = 2832.9375
$ws.Cells.Item(113, 11).Value = 2832.9375
$ws.Cells.Item(113, 13).Value = -662.9375

$ws.Cells.Item(122, 8).Value = 3004559.2
$ws.Cells.Item(122, 9).Value = 5003952
$ws.Cells.Item(122, 11).Value = 15011856
$ws.Cells.Item(122, 13).Value = -15009406

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(97, 8).Value = 10280
$ws.Cells.Item(97, 10).Value = 10280
$ws.Cells.Item(97, 12).Value = 10280
$ws.Cells.Item(97, 14).Value = -12262

$ws.Cells.Item(132, 8).Value = 4235.5
$ws.Cells.Item(132, 9).Value = 3730.8
$ws.Cells.Item(132, 10).Value = 4977.706
$ws.Cells.Item(132, 11).Value = 11192.4
$ws.Cells.Item(132, 12).Value = 14933.118
$ws.Cells.Item(132, 13).Value = -8662.400000000001
$ws.Cells.Item(132, 14).Value = -19993.118

$ws.Cells.Item(136, 8).Value = 2249.0605
$ws.Cells.Item(136, 9).Value = 2171
$ws.Cells.Item(136, 11).Value = 6513
$ws.Cells.Item(136, 13).Value = -3963
